$d = $word.ActiveDocument

$pkgHeader = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgFooter = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$newParas = '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Possibilit' + [char]0xE9 + ' d' + [char]0x2019 + 'ajouter un jeu si l' + [char]0x2019 + 'utilisateur est administrateur</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Possibilit' + [char]0xE9 + ' de supprimer un jeu si l' + [char]0x2019 + 'utilisateur est administrateur</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xml = $pkgHeader + $newParas + $pkgFooter

$p25 = $d.Paragraphs.Item(25)
$target = $d.Range($p25.Range.Start, $p25.Range.Start)
$target.InsertXML($xml)

Write-Host "done"
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($i -ge 22 -and $i -le 32) {
        Write-Host "$i : [$($p.Range.Text)]"
    }
    $i = $i + 1
}
